$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.466.77"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.793.68"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.67"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.61"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.296"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.053.59"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.784.84"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.636"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.449.36"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.27"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.85"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.25"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0800"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.19"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.17"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.06"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.25"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.50"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.81"
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.23"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0521"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.88"
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.433.09"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.667"
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.61"
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.47"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0526"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.11"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.949.98"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.67"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("E51").Value = "  +0.07%  "

Write-Output "Updated cryptos list"
